# Reposition/resize the existing chart's anchor on "Hoja1" so that it spans
# from (col C, row 14) to just past (col K, row 30) instead of its previous
# anchor from (col A, row 14) to (col J, row 29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$co = $ws.ChartObjects().Item(1)

# Target geometry (in points) computed from the new two-cell anchor:
#   from: col=2 (0-based) colOff=0 EMU, row=13 (0-based) rowOff=42861 EMU
#   to:   col=10 (0-based) colOff=323850 EMU, row=29 (0-based) rowOff=9524 EMU
$co.Left = 145.4287109375
$co.Top = 198.37488188976377
$co.Width = 493
$co.Height = 237.37503937007875
